$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" -> "Impact" bullet list from six
# job-duty-style bullets down to four impact-focused accomplishment bullets.
#
# Some of the original bullet phrasing is duplicated verbatim elsewhere in
# the resume (Professional Experience section), so every lookup below is
# scoped to paragraph indices inside the Key Achievements block rather than
# using an unscoped document-wide Find/Replace.
# ---------------------------------------------------------------------------

function Find-ParaIndex($searchText, $fromIdx, $toIdx) {
    for ($idx = $fromIdx; $idx -le $toIdx; $idx++) {
        $p = $d.Paragraphs.Item($idx)
        if ($p.Range.Text -like "*$searchText*") {
            return $idx
        }
    }
    return -1
}

# Locate the section heading.
$headingIdx = Find-ParaIndex "KEY ACHIEVEMENTS AND IMPACT" 1 $d.Paragraphs.Count
$scanEnd = [Math]::Min($headingIdx + 15, $d.Paragraphs.Count)

# 1) Delete the two bullets being removed entirely: "Discovered systematic
#    race coding errors..." and "Achieved 87% prediction accuracy...".
$idxRace = Find-ParaIndex "Discovered systematic race coding errors" $headingIdx $scanEnd
$idx87 = Find-ParaIndex "Achieved 87% prediction accuracy" $headingIdx $scanEnd

$pRace = $d.Paragraphs.Item($idxRace)
$p87 = $d.Paragraphs.Item($idx87)
$killRange = $d.Range($pRace.Range.Start, $p87.Range.End)
$killRange.Delete()

# 2) Rewrite the remaining four bullets with impact-focused accomplishment
#    text (each lookup re-scanned/re-scoped after the deletion above).
$scanEnd = [Math]::Min($headingIdx + 15, $d.Paragraphs.Count)

$idxBuilt = Find-ParaIndex "Built redistricting platform used by thousands of analysts nationwide" $headingIdx $scanEnd
$pBuilt = $d.Paragraphs.Item($idxBuilt)
$pBuilt.Range.Text = "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"

$idxEtl = Find-ParaIndex "Designed ETL pipelines using PySpark" $headingIdx $scanEnd
$pEtl = $d.Paragraphs.Item($idxEtl)
$pEtl.Range.Text = "• `$4.7M savings enabled nonprofit access"

$idxTrig = Find-ParaIndex "Trigonometric algorithm for boundary estimation" $headingIdx $scanEnd
$pTrig = $d.Paragraphs.Item($idxTrig)
$pTrig.Range.Text = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"

$idxCloud = Find-ParaIndex "Built cloud-based data warehouse solutions" $headingIdx $scanEnd
$pCloud = $d.Paragraphs.Item($idxCloud)
$pCloud.Range.Text = "• Real-time collaboration at national scale"

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
